$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-07"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 04-07)"

# Update the April (row 5) and Total (row 14) figures for the "Total" column (I)
$ws.Range("I5").Value = 24
$ws.Range("I14").Value = 458
